# Team_PER_2017: fix the PER (Player Efficiency Rating) bug — the sheet
# previously held raw per-team point totals in column C paired against the
# wrong team names; the fix reshuffles which team sits on each row and
# replaces column C with the corrected (much smaller) per-game PER-style
# ratios.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Team (column B) order for rows 2..31 (row 1 is the "Team"/"PER" header).
$teams = @(
    "POR", "CLE", "DAL", "MIA", "OKC", "ATL", "WAS", "MIL", "LAC", "SAS",
    "DET", "ORL", "UTA", "MEM", "HOU", "NOP", "DEN", "LAL", "GSW", "IND",
    "CHO", "CHI", "PHI", "BOS", "BRK", "TOR", "SAC", "PHO", "NYK", "MIN"
)

# Corrected PER values (column C) for rows 2..31, in the same row order.
$values = @(
    12.54615384615385, 13.74285714285714, 10.46875,           12.99285714285714,
    13.91666666666667, 10.32666666666667, 13.22142857142857,  12.92666666666667,
    13.48,              13.99375,          13.06,              9.456250000000001,
    14.48666666666667,  12.85333333333333, 14.01538461538462,  13.3,
    15.15333333333333,  13.04285714285714, 15.43333333333333,  12.92,
    12.19285714285714,  12.1923076923077,  13.5,               14.78666666666667,
    12.66428571428571,  14.44285714285714, 13.32307692307692,  12.25882352941176,
    13.64285714285714,  13.15
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $teams[$i]
    $ws.Cells.Item($r, 3).Value = $values[$i]
}
